$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "ID" column (A2:A12) values from O_0xx to MAT_0xx
for ($i = 2; $i -le 12; $i++) {
    $n = $i - 1
    $ws.Cells.Item($i, 1).Value = "MAT_{0:D3}" -f $n
}

# Update the active selection on the sheet
$ws.Range("B15").Select()
